# LeetCode Practice Tracker update
# Adds new Queue-topic problems (rows 125-131), relocates "Valid Anagram"
# (String topic) that used to live at row 127 down to row 142, and marks
# three problems as solved per the commit message:
#   1. Design Circular Queue
#   2. Implement Queue using Stacks
#   3. Time Needed to Buy Tickets

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed the brand-new shared strings first, in the same order they end up in
# the workbook's string table, so the resulting sharedStrings.xml lines up
# with the source edit.
$ws.Range("C126").Value = "Design Circular Queue"
$ws.Range("C130").Value = "Design Hit Counter"
$ws.Range("G125").Value = "O(n) Amortized O(1)"
$ws.Range("C127").Value = "Moving Average from Data Stream"
$ws.Range("C128").Value = "Number of Recent Calls"
$ws.Range("C131").Value = "Time Needed to Buy Tickets"
$ws.Range("C129").Value = "Dota2 Senate"
$ws.Range("H126").Value = "O(capacity)"

# --- Row 125: Implement Queue using Stacks (solved) ---
$ws.Range("B125").Value = "Queue"
$ws.Range("C125").Value = "Implement Queue using Stacks"
$ws.Range("D125").Value = "Easy"
$ws.Range("E125").Value = "Done"
$ws.Range("H125").Value = "O(n)"

# --- Row 126: Design Circular Queue (solved) ---
$ws.Range("B126").Value = "Queue"
$ws.Range("D126").Value = "Medium"
$ws.Range("E126").Value = "Done"
$ws.Range("G126").Value = "O(1)"

# --- Row 127: Moving Average from Data Stream (not yet solved) ---
# Previously held "Valid Anagram" (String) - that data moves to row 142.
$ws.Range("B127").Value = "Queue"
$ws.Range("D127").Clear()
$ws.Range("E127").Clear()

# --- Row 128: Number of Recent Calls (not yet solved) ---
# Previously held "Implement Queue using Stacks" placeholder entry, now moved
# (solved) up to row 125, so this row is repurposed.
$ws.Range("D128").Clear()
$ws.Range("E128").Clear()

# --- Row 129: Dota2 Senate (not yet solved) ---
$ws.Range("B129").Value = "Queue"

# --- Row 130: Design Hit Counter (not yet solved) ---
$ws.Range("B130").Value = "Queue"

# --- Row 131: Time Needed to Buy Tickets (solved, new row) ---
$ws.Range("A131").Value = 131
$ws.Range("B131").Value = "Queue"
$ws.Range("D131").Value = "Easy"
$ws.Range("E131").Value = "Done"
$ws.Range("G131").Value = "O(n)"
$ws.Range("H131").Value = "O(1)"

# --- Row 142: relocated "Valid Anagram" (String) entry, no ID value ---
$ws.Range("B142").Value = "String"
$ws.Range("C142").Value = "Valid Anagram"
$ws.Range("D142").Value = "Easy"
$ws.Range("E142").Value = "To Do"

# --- Date Solved (column F) cells: copy the date number format from an
#     existing date cell so the style matches (s="3"), then set the value.
$ws.Range("F124").Copy()
$ws.Range("F125").PasteSpecial(-4122)
$ws.Range("F126").PasteSpecial(-4122)
$ws.Range("F131").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$ws.Range("F125").Value = 46037
$ws.Range("F126").Value = 46037
$ws.Range("F131").Value = 45672

# --- View state: selection moves to C128 (matches the diff) ---
$ws.Range("C128").Select()
